$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = 45574
$ws.Range("C2").Value = 45573
$ws.Range("D2").Value = 0.49
$ws.Range("E2").Value = 147024010
$ws.Range("F2").Value = 147000000

# Data for new rows 3-11
$rows = @(
    @{ Row=3;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=2.1;      E=2000;           F=76027595;   G=0; H="BCHIEJ0717"; I="VENTA";  J="RENTA FIJA" },
    @{ Row=4;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=2.51;     E=2000;           F=80139302;   G=0; H="BESTJ41008"; I="VENTA";  J="RENTA FIJA" },
    @{ Row=5;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=2.21;     E=15500;          F=597519830;  G=0; H="BBBVK61113"; I="VENTA";  J="RENTA FIJA" },
    @{ Row=6;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=2.51;     E=5000;           F=200348256;  G=0; H="BESTJ41008"; I="VENTA";  J="RENTA FIJA" },
    @{ Row=7;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=0;        E=33025154.912;   F=33012410;   G=0; H="SOCOVESA";   I="VENTA";  J="SIMULTANEA" },
    @{ Row=8;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=0;        E=390210469.11;   F=390000887;  G=0; H="BESALCO";    I="VENTA";  J="SIMULTANEA" },
    @{ Row=9;  A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=0;        E=30621570;       F=30471350;   G=0; H="LTM";        I="VENTA";  J="SIMULTANEA" },
    @{ Row=10; A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=16532.39; E=2773;           F=45844317;   G=0; H="CFINHRFLA";  I="VENTA";  J="RENTA VARIABLE" },
    @{ Row=11; A="FONDO DE INVERSION NEVASA AHORRO"; B=45573; C=45573; D=15976.83; E=557;            F=8899095;    G=0; H="CFINHRFLB";  I="COMPRA"; J="RENTA VARIABLE" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
}
